# fix in recode.others function
#
# The three data rows describing the "a7_site_management_receive_compensation"
# group (rows 4-6) get collapsed into a single merged row (row 4), which
# shifts every row below up by two. After that shift, the issue note on the
# two "humanitarian_agency" rows (new rows 5 and 6) is updated to flag the
# extra ";)" note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collapse old rows 4+5+6 into a single row by deleting rows 6 and 5
# (bottom-up so row indices stay valid while deleting).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(5).Delete()

# Row 4 now holds the old "a7_site_management_receive_compensation" data;
# rewrite it to the merged "myother" / phone-number row.
$ws.Cells.Item(4, 4).Value = "a7_1_site_management_receive_compensation_other"
$ws.Cells.Item(4, 5).Value = "myother;Phone Number"
$ws.Cells.Item(4, 6).Value = "myother"
$ws.Cells.Item(4, 7).Value = "Recode elsewhere other"

# The humanitarian_agency rows (old rows 7 and 8) are now rows 5 and 6;
# only the issue note changes.
$ws.Cells.Item(5, 7).Value = "Recode elsewhere other;)"
$ws.Cells.Item(6, 7).Value = "Recode elsewhere other;)"

Write-Output "done"
